$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Brand Bulk Upload template polish:
#  - remove the "description" column and its sample data ("Brand Description")
#  - remove the unused "Active"/"Inactive" text labels, replacing status /
#    featured sample values with the numeric flags 1 / 0
#  - move "logo" to the last column (after the new blank "components" column)
# ---------------------------------------------------------------------------

# Step 1: drop the trailing (empty) "components" column - it is the very last
# column so this is a clean delete that does not disturb any other widths.
$ws.Columns("J").Delete()

# Step 2: drop the "description" column (was D). Everything from old E..I
# shifts left into D..H.
$ws.Columns("D").Delete()

# Layout is now:
#   A name | B slug | C logo | D meta_title | E meta_keyword |
#   F meta_description | G status | H featured
# Remember logo's current header/value so it can be re-created at the end.
$logoHeader = $ws.Range("C1").Value2
$logoValue  = $ws.Range("C2").Value2

# Step 3: shift meta_title..featured one column to the left (into the old
# logo column and onward), turning status/featured into numeric 1/0 flags.
$ws.Range("C1").Value = $ws.Range("D1").Value2
$ws.Range("C2").Value = $ws.Range("D2").Value2

$ws.Range("D1").Value = $ws.Range("E1").Value2
$ws.Range("D2").Value = $ws.Range("E2").Value2

$ws.Range("E1").Value = $ws.Range("F1").Value2
$ws.Range("E2").Value = $ws.Range("F2").Value2

$ws.Range("F1").Value = $ws.Range("G1").Value2
$ws.Range("F2").Value = 1

$ws.Range("G1").Value = $ws.Range("H1").Value2
$ws.Range("G2").Value = 0

# Step 4: new blank "components" header in H, matching the existing bold
# header style (reuse font attributes instead of Range.Style so the engine
# maps it back onto the very same cellXf instead of minting a new one).
$ws.Range("H1").Value = "components"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").Font.Name = "Arial"
$ws.Range("H1").Font.Size = 10
$ws.Range("H1").Font.Color = $ws.Range("A1").Font.Color
$ws.Range("H2").ClearContents()

# Step 5: logo becomes the last (I) column, keeping its bold header style.
$ws.Range("I1").Value = $logoHeader
$ws.Range("I1").Font.Bold = $true
$ws.Range("I1").Font.Name = "Arial"
$ws.Range("I1").Font.Size = 10
$ws.Range("I1").Font.Color = $ws.Range("A1").Font.Color
$ws.Range("I2").Value = $logoValue

# Restore the selection to E2 (matches the refreshed sheet view).
$ws.Range("E2").Select()
